$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20 (ID 19 / File 019): refresh the "2 minute" reminder text to a
#     generic version that no longer hard-codes specific speeds ---
$ws.Range("C20").Value = "Bạn đã chạy được 2 phút với tốc độ hiện tại, hãy bấm phím speed cộng để tăng hoặc speed trừ để giảm tốc độ."

# --- Row 21 (ID 20 / File 020): now holds the "Volume sound" clip ---
$ws.Range("C21").Value = "Volume sound"

# --- Rows 24-26 and 28: the old per-speed "2 minute" reminder rows and the
#     "Hãy nhấn phím Stop..." row are retired. Column A/B content is
#     cleared (formatting kept) while column C is fully cleared (format
#     reverts to the column default, so the cell drops out entirely) ---
$ws.Range("A24:B26").ClearContents()
$ws.Range("C24:C26").Clear()

$ws.Range("A28:B28").ClearContents()
$ws.Range("C28").Clear()

# --- Row 27: becomes the new "Sleep mode 10p" row - A/B cleared, C set ---
$ws.Range("A27:B27").ClearContents()
$ws.Range("C27").Value = "Sleep mode 10p"

# --- Row 22 (ID 21 / File 021): "speed up" voice prompt ---
$ws.Range("C22").Value = "Tăng tốc độ, chú ý giữ an toàn nhé"

# --- Row 23 (File 022): "speed down" voice prompt. The numeric ID in
#     column A is retired (cleared) while File Name/Content remain ---
$ws.Range("A23").ClearContents()
$ws.Range("C23").Value = "Giảm tốc độ, chú ý giữ an toàn nhé"

# --- Row 29: entirely retired (old "028 / Volume sound" row) ---
$ws.Range("A29:C29").Clear()

# --- Match the author's final selection state ---
$ws.Range("C22").Select()
